$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2: replace old prompt text with new greeting
$ws.Range("C2").Value = "Dobar dan! ☀️ Kako vam mogu pomoći?"

# Adjust row heights for rows 3 and 4
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 45

# Add new row 5 data
$ws.Range("B5").Value = "kartice"
$ws.Range("A5").Value = "HR523"
$ws.Range("C5").Value = "Kod nas možete ugovoriti dvije vrsta kartica ovisno o tome što vam treba.
Debitne kartice koje su vezane uz tekući, žiro, zaštićeni ili multivalutni tekući račun. Kada plaćate debitnim karticama, novac se odmah skida s računa.
Kreditne kartice kod kojih ne plaćate kupnju odmah, nego kasnije.
Koje vas kartice zanimaju?
Debitne kartice
Kreditne kartice"
$ws.Rows.Item(5).RowHeight = 180
$ws.Range("C5").WrapText = $true

# Adjust column width for column C (closest achievable to target 37.42578125)
$ws.Columns.Item(3).ColumnWidth = 36.7

# Update selection
$ws.Range("C5").Select() | Out-Null
